$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.901.66'
$ws.Range("E2").Value = '  +4.93%  '
$ws.Range("D3").Value = '2.287.50'
$ws.Range("E3").Value = '  +2.40%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '231.84'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.43%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.629'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.67%  '
$ws.Range("E7").Value = '  +5.31%  '
$ws.Range("E8").Value = '  +0.10%  '
$ws.Range("E9").Value = '  +4.51%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0953'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.93%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '57.69'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.62%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '26.39'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +15.64%  '
$ws.Range("E13").Value = '  +0.20%  '
$ws.Range("D14").Value = '2.627.32'
$ws.Range("E14").Value = '  +2.42%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.81'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.09%  '
$ws.Range("E16").Value = '  +5.36%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.816'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.15%  '
$ws.Range("D18").Value = '2.274.41'
$ws.Range("E18").Value = '  +1.11%  '
$ws.Range("D19").Value = '43.790.71'
$ws.Range("E19").Value = '  +4.86%  '
$ws.Range("D20").Value = '0.0₃0946'
$ws.Range("E20").Value = '  +4.07%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '73.35'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.95%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.18'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.83%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '250.77'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.87%  '
$ws.Range("E24").Value = '  -0.01%  '
$ws.Range("E25").Value = '  +8.40%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.27'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.22%  '
$ws.Range("E27").Value = '  +1.05%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '171.61'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.07%  '
$ws.Range("E29").Value = '  -2.96%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '20.56'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.88%  '
$ws.Range("E31").Value = '  +2.14%  '
$ws.Range("E33").Value = '  +0.02%  '
$ws.Range("E34").Value = '  +7.07%  '
$ws.Range("E35").Value = '  +1.68%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.73'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.14%  '
$ws.Range("E37").Value = '  +0.49%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.71'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.50%  '
$ws.Range("E39").Value = '  -1.76%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0249'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.07%  '
$ws.Range("E41").Value = '  +0.22%  '
$ws.Range("B42").Value = 'FTXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.74'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +6.87%  '
$ws.Range("B43").Value = 'Celestia'
$ws.Range("C43").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.01'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +26.95%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.52'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.15%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.000221'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -7.93%  '
$ws.Range("E46").Value = '  -0.46%  '
$ws.Range("E47").Value = '  +0.76%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '98.10'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.96%  '
$ws.Range("D49").Value = '1.488.69'
$ws.Range("E49").Value = '  +1.11%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '16.89'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.14%  '
$ws.Range("E51").Value = '  +3.47%  '
